$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.095.81'
$ws.Range('E2').Value = '  -2.01%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.806.28'
$ws.Range('E3').Value = '  +0.02%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  -0.16%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.60'
$ws.Range('E5').Value = '  -0.41%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9994'
$ws.Range('E6').Value = '  -0.17%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5307'
$ws.Range('E7').Value = '  -2.01%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3972'
$ws.Range('E8').Value = '  +4.73%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07574'
$ws.Range('E9').Value = '  +0.79%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.58'
$ws.Range('E10').Value = '  -1.98%  '

$ws.Range('E11').Value = '  -2.13%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9992'
$ws.Range('E12').Value = '  -0.24%  '

$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.255'
$ws.Range('E13').Value = '  +1.61%  '

$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.577'
$ws.Range('E14').Value = '  +2.90%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.44'
$ws.Range('E15').Value = '  -1.00%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.812.10'
$ws.Range('E16').Value = '  +0.74%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '89.10'
$ws.Range('E17').Value = '  -1.22%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001065'
$ws.Range('E18').Value = '  -0.09%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06568'
$ws.Range('E19').Value = '  +0.98%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9996'
$ws.Range('E20').Value = '  -0.11%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.31'
$ws.Range('E21').Value = '  -0.37%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.017'
$ws.Range('E22').Value = '  +1.01%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.109.43'
$ws.Range('E23').Value = '  -1.98%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.12'
$ws.Range('E24').Value = '  -0.22%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.088'
$ws.Range('E25').Value = '  +0.63%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '156.26'
$ws.Range('E26').Value = '  -3.07%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.30'
$ws.Range('E27').Value = '  -0.97%  '

$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.377'
$ws.Range('E28').Value = '  +1.62%  '

$ws.Range('B29').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C29').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.009.90'
$ws.Range('E29').Value = '  +0.10%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '122.93'
$ws.Range('E30').Value = '  -0.17%  '

$ws.Range('E31').Value = '  +3.57%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.095'
$ws.Range('E32').Value = '  -4.65%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.665'
$ws.Range('E33').Value = '  -0.52%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.548'
$ws.Range('E34').Value = '  -2.07%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07079'
$ws.Range('E35').Value = '  +6.40%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.2216'
$ws.Range('E36').Value = '  -2.45%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.178'
$ws.Range('E37').Value = '  +2.86%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02294'
$ws.Range('E38').Value = '  -0.74%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.490'
$ws.Range('E39').Value = '  -1.47%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '11.24'
$ws.Range('E40').Value = '  -0.07%  '

$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.194'
$ws.Range('E41').Value = '  -0.59%  '

$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6167'
$ws.Range('E42').Value = '  -1.16%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.406'
$ws.Range('E43').Value = '  -3.12%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.42'
$ws.Range('E44').Value = '  +0.85%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.681'
$ws.Range('E45').Value = '  -0.63%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5747'
$ws.Range('E46').Value = '  -1.87%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '125.03'
$ws.Range('E47').Value = '  -1.37%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.934'
$ws.Range('E48').Value = '  -1.32%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.186'
$ws.Range('E49').Value = '  -0.12%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06812'
$ws.Range('E50').Value = '  -1.41%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '71.34'
$ws.Range('E51').Value = '  -1.96%  '
